$wb = $excel.ActiveWorkbook

# --- item_level sheet: rename two VAR_NAMES values ---
$ws1 = $wb.Worksheets.Item("item_level")
$ws1.Range("A13").Value = "date_of_data_provision"
$ws1.Range("A14").Value = "date_of_archiving"

# --- checks sheet: update CONTRADICTION_TERM formulas to reference the new names ---
$ws2 = $wb.Worksheets.Item("checks")
$ws2.Range("C3").Value = '[date_of_submission] =  "" and [date_of_data provision] <> ""'
$ws2.Range("C4").Value = '[date_of_submission] =  "" and [date_of_archiving] <> ""'
$ws2.Range("C5").Value = '[date_of_data provision] =  "" and [date_of_archiving] <> ""'

# --- sheet view / selection state changes ---
$ws2.Activate()
$ws2.Range("C8").Select()

$ws1.Activate()
$ws1.Range("B16").Select()
